$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 7461.143
$ws.Range("J40").Value = 8068
$ws.Range("L40").Value = 8068
$ws.Range("N40").Value = -8418
$ws.Range("H70").Value = 2912.8572
$ws.Range("I70").Value = 2156
$ws.Range("J70").Value = 3480.5
$ws.Range("K70").Value = 6468
$ws.Range("L70").Value = 10441.5
$ws.Range("M70").Value = -6198
$ws.Range("N70").Value = -10981.5
$ws.Range("H73").Value = 2912.8572
$ws.Range("I73").Value = 2156
$ws.Range("J73").Value = 3480.5
$ws.Range("K73").Value = 6468
$ws.Range("L73").Value = 10441.5
$ws.Range("M73").Value = -5532
$ws.Range("N73").Value = -12313.5
$ws.Range("H100").Value = 2747.25
$ws.Range("I100").Value = 2747.25
$ws.Range("K100").Value = 2747.25
$ws.Range("M100").Value = -2206.25
$ws.Range("H103").Value = 2613.4
$ws.Range("I103").Value = 3032.3333
$ws.Range("K103").Value = 9096.999899999999
$ws.Range("M103").Value = -8510.999899999999
$ws.Range("H128").Value = 150000
$ws.Range("J128").Value = 150000
$ws.Range("L128").Value = 150000
$ws.Range("N128").Value = -159960
$ws.Range("H133").Value = 93991.75
$ws.Range("J133").Value = 93991.75
$ws.Range("L133").Value = 93991.75
$ws.Range("N133").Value = -104111.75
$ws.Range("H137").Value = 4296
$ws.Range("I137").Value = 0
$ws.Range("K137").Value = 0
$ws.Range("M137").ClearContents()
$ws.Range("H138").Value = 309218.4
$ws.Range("I138").Value = 9312
$ws.Range("J138").Value = 341208.44
$ws.Range("K138").Value = 27936
$ws.Range("L138").Value = 1023625.32
$ws.Range("M138").Value = -22796
$ws.Range("N138").Value = -1033905.32
$ws.Range("H141").Value = 3350
$ws.Range("I141").Value = 3180.875
$ws.Range("J141").Value = 3620.6
$ws.Range("K141").Value = 9542.625
$ws.Range("L141").Value = 10861.8
$ws.Range("M141").Value = -4362.625
$ws.Range("N141").Value = -21221.8

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2249.68
$ws.Range("I2").Value = 2243.4167
$ws.Range("K2").Value = 2243.4167
$ws.Range("M2").Value = -2130.4167
$ws.Range("H54").Value = 15000
$ws.Range("J54").Value = 15000
$ws.Range("L54").Value = 15000
$ws.Range("N54").Value = -16538
$ws.Range("H74").Value = 176109.94
$ws.Range("I74").Value = 224244.16
$ws.Range("K74").Value = 224244.16
$ws.Range("M74").Value = -223370.16
$ws.Range("H77").Value = 176109.94
$ws.Range("I77").Value = 224244.16
$ws.Range("K77").Value = 1121220.8
$ws.Range("M77").Value = -1116852.8
$ws.Range("H80").Value = 74805
$ws.Range("J80").Value = 74805
$ws.Range("L80").Value = 74805
$ws.Range("N80").Value = -76801
$ws.Range("H83").Value = 74805
$ws.Range("J83").Value = 74805
$ws.Range("L83").Value = 224415
$ws.Range("N83").Value = -234399
$ws.Range("H102").Value = 5490.9443
$ws.Range("I102").Value = 5425.6333
$ws.Range("K102").Value = 5425.6333
$ws.Range("M102").Value = -3803.6333
$ws.Range("H116").Value = 2249.68
$ws.Range("I116").Value = 2243.4167
$ws.Range("K116").Value = 2243.4167
$ws.Range("M116").Value = 50.58329999999978

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2249.68
$ws.Range("I3").Value = 2243.4167
$ws.Range("K3").Value = 2243.4167
$ws.Range("M3").Value = -2129.4167
$ws.Range("H57").Value = 30285.428
$ws.Range("J57").Value = 30285.428
$ws.Range("L57").Value = 30285.428
$ws.Range("N57").Value = -31725.428
$ws.Range("H80").Value = 83333620
$ws.Range("J80").Value = 476
$ws.Range("L80").Value = 476
$ws.Range("N80").Value = -2472
$ws.Range("H83").Value = 83333620
$ws.Range("J83").Value = 476
$ws.Range("L83").Value = 2380
$ws.Range("N83").Value = -12364
$ws.Range("H136").Value = 30285.428
$ws.Range("J136").Value = 30285.428
$ws.Range("L136").Value = 30285.428
$ws.Range("N136").Value = -40485.428

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3601.9656
$ws.Range("I31").Value = 2739.3
$ws.Range("J31").Value = 4056
$ws.Range("K31").Value = 2739.3
$ws.Range("L31").Value = 4056
$ws.Range("M31").Value = -2444.3
$ws.Range("N31").Value = -4646
$ws.Range("H34").Value = 3601.9656
$ws.Range("I34").Value = 2739.3
$ws.Range("J34").Value = 4056
$ws.Range("K34").Value = 2739.3
$ws.Range("L34").Value = 4056
$ws.Range("M34").Value = -2537.3
$ws.Range("N34").Value = -4460
$ws.Range("H58").Value = 4109.5
$ws.Range("I58").Value = 3783.8462
$ws.Range("K58").Value = 3783.8462
$ws.Range("M58").Value = -3580.8462
$ws.Range("H99").Value = 5732.3
$ws.Range("I99").Value = 6003.4287
$ws.Range("K99").Value = 6003.4287
$ws.Range("M99").Value = -4505.4287
$ws.Range("H126").Value = 5732.3
$ws.Range("I126").Value = 6003.4287
$ws.Range("K126").Value = 18010.2861
$ws.Range("M126").Value = -15540.2861
$ws.Range("H132").Value = 4064.1667
$ws.Range("I132").Value = 3835.756
$ws.Range("J132").Value = 4784.5386
$ws.Range("K132").Value = 11507.268
$ws.Range("L132").Value = 14353.6158
$ws.Range("M132").Value = -8977.268
$ws.Range("N132").Value = -19413.6158
$ws.Range("H134").Value = 3477.2
$ws.Range("I134").Value = 2512.68
$ws.Range("J134").Value = 8299.799999999999
$ws.Range("K134").Value = 7538.039999999999
$ws.Range("L134").Value = 24899.4
$ws.Range("M134").Value = -5003.039999999999
$ws.Range("N134").Value = -29969.4
$ws.Range("H136").Value = 4109.5
$ws.Range("I136").Value = 3783.8462
$ws.Range("K136").Value = 11351.5386
$ws.Range("M136").Value = -8801.5386
$ws.Range("H141").Value = 822496
$ws.Range("J141").Value = 822496
$ws.Range("L141").Value = 822496
$ws.Range("N141").Value = -832856

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 631.5
$ws.Range("I38").Value = 317
$ws.Range("K38").Value = 951
$ws.Range("M38").Value = -604
$ws.Range("H58").Value = 2174.5
$ws.Range("J58").Value = 3999
$ws.Range("L58").Value = 11997
$ws.Range("N58").Value = -12253
$ws.Range("H80").Value = 3749.5
$ws.Range("J80").Value = 4332.6665
$ws.Range("L80").Value = 12997.9995
$ws.Range("N80").Value = -14869.9995
$ws.Range("H83").Value = 3749.5
$ws.Range("J83").Value = 4332.6665
$ws.Range("L83").Value = 38993.9985
$ws.Range("N83").Value = -48353.9985
$ws.Range("H113").Value = 5480
$ws.Range("J113").Value = 5777.778
$ws.Range("L113").Value = 17333.334
$ws.Range("N113").Value = -21673.334

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 747.7143
$ws.Range("I2").Value = 879.8
$ws.Range("J2").Value = 417.5
$ws.Range("K2").Value = 879.8
$ws.Range("L2").Value = 417.5
$ws.Range("M2").Value = -766.8
$ws.Range("N2").Value = -643.5
$ws.Range("H32").Value = 49996.332
$ws.Range("J32").Value = 49996.332
$ws.Range("L32").Value = 49996.332
$ws.Range("N32").Value = -50588.332
$ws.Range("H132").Value = 2802.8276
$ws.Range("I132").Value = 1958.5
$ws.Range("K132").Value = 5875.5
$ws.Range("M132").Value = -3345.5

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2294.1765
$ws.Range("J46").Value = 2020.8334
$ws.Range("L46").Value = 2020.8334
$ws.Range("N46").Value = -2396.8334
$ws.Range("H94").Value = 587625
$ws.Range("J94").Value = 587625
$ws.Range("L94").Value = 587625
$ws.Range("N94").Value = -588977
$ws.Range("H136").Value = 4452.778
$ws.Range("I136").Value = 3100.0715
$ws.Range("J136").Value = 9187.25
$ws.Range("K136").Value = 9300.2145
$ws.Range("L136").Value = 27561.75
$ws.Range("M136").Value = -6750.2145
$ws.Range("N136").Value = -32661.75

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 6443.275
$ws.Range("J132").Value = 8285.571
$ws.Range("L132").Value = 24856.713
$ws.Range("N132").Value = -29916.713
$ws.Range("H136").Value = 2272.8462
$ws.Range("I136").Value = 1582.4117
$ws.Range("K136").Value = 4747.2351
$ws.Range("M136").Value = -2197.2351
